$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 13386.223
$ws.Range("I70").Value = 993.5
$ws.Range("K70").Value = 2980.5
$ws.Range("M70").Value = -2710.5
$ws.Range("H73").Value = 13386.223
$ws.Range("I73").Value = 993.5
$ws.Range("K73").Value = 2980.5
$ws.Range("M73").Value = -2044.5
$ws.Range("H86").Value = 4466.6665
$ws.Range("I86").Value = 900
$ws.Range("J86").Value = 6250
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 6250
$ws.Range("M86").Value = 223
$ws.Range("N86").Value = -8496
$ws.Range("H89").Value = 4466.6665
$ws.Range("I89").Value = 900
$ws.Range("J89").Value = 6250
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 31250
$ws.Range("M89").Value = 1116
$ws.Range("N89").Value = -42482
$ws.Range("H98").Value = 3380.5715
$ws.Range("I98").Value = 3650.6667
$ws.Range("K98").Value = 3650.6667
$ws.Range("M98").Value = -2152.6667
$ws.Range("H122").Value = 3380.5715
$ws.Range("I122").Value = 3650.6667
$ws.Range("K122").Value = 10952.0001
$ws.Range("M122").Value = -8502.000100000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3365.3333
$ws.Range("I61").Value = 3162
$ws.Range("K61").Value = 3162
$ws.Range("M61").Value = -2950
$ws.Range("H122").Value = 9746.5
$ws.Range("I122").Value = 9996.111000000001
$ws.Range("K122").Value = 29988.333
$ws.Range("M122").Value = -27538.333
$ws.Range("H136").Value = 3365.3333
$ws.Range("I136").Value = 3162
$ws.Range("K136").Value = 9486
$ws.Range("M136").Value = -6936

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2863.2
$ws.Range("I20").Value = 2172
$ws.Range("J20").Value = 3900
$ws.Range("K20").Value = 2172
$ws.Range("L20").Value = 3900
$ws.Range("M20").Value = -1925
$ws.Range("N20").Value = -4394

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 502000.75
$ws.Range("I3").Value = 1002500
$ws.Range("J3").Value = 1501.5
$ws.Range("K3").Value = 1002500
$ws.Range("L3").Value = 1501.5
$ws.Range("M3").Value = -1002387
$ws.Range("N3").Value = -1727.5
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 484.2857
$ws.Range("I11").Value = 481.66666
$ws.Range("K11").Value = 1444.99998
$ws.Range("M11").Value = -1304.99998
$ws.Range("H92").Value = 2275
$ws.Range("J92").Value = 3500
$ws.Range("L92").Value = 10500
$ws.Range("N92").Value = -12996
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5250
$ws.Range("I7").Value = 2500
$ws.Range("K7").Value = 2500
$ws.Range("M7").Value = -2388
$ws.Range("H8").Value = 5250
$ws.Range("I8").Value = 2500
$ws.Range("K8").Value = 2500
$ws.Range("M8").Value = -2361
$ws.Range("H20").Value = 15555
$ws.Range("I20").Value = 15555
$ws.Range("K20").Value = 15555
$ws.Range("M20").Value = -15310
$ws.Range("H80").Value = 4459.8
$ws.Range("I80").Value = 3750
$ws.Range("J80").Value = 4933
$ws.Range("K80").Value = 3750
$ws.Range("L80").Value = 4933
$ws.Range("M80").Value = -2752
$ws.Range("N80").Value = -6929
$ws.Range("H83").Value = 4459.8
$ws.Range("I83").Value = 3750
$ws.Range("J83").Value = 4933
$ws.Range("K83").Value = 18750
$ws.Range("L83").Value = 24665
$ws.Range("M83").Value = -13758
$ws.Range("N83").Value = -34649
$ws.Range("H102").Value = 2524.75
$ws.Range("I102").Value = 2524.75
$ws.Range("K102").Value = 2524.75
$ws.Range("M102").Value = -902.75
$ws.Range("H122").Value = 1187.25
$ws.Range("I122").Value = 1133
$ws.Range("K122").Value = 3399
$ws.Range("M122").Value = -949

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4336.875
$ws.Range("I7").Value = 3924.5
$ws.Range("K7").Value = 3924.5
$ws.Range("M7").Value = -3812.5
$ws.Range("H16").Value = 2500
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 2500
$ws.Range("M16").Value = -2330
$ws.Range("H21").Value = 6811.25
$ws.Range("I21").Value = 5748.3335
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 5748.3335
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -5574.3335
$ws.Range("N21").Value = -10348
$ws.Range("H22").Value = 1850.3334
$ws.Range("I22").Value = 1850.3334
$ws.Range("K22").Value = 1850.3334
$ws.Range("M22").Value = -1555.3334
$ws.Range("H27").Value = 1850.3334
$ws.Range("I27").Value = 1850.3334
$ws.Range("K27").Value = 1850.3334
$ws.Range("M27").Value = -1743.3334
$ws.Range("H46").Value = 3944.889
$ws.Range("I46").Value = 2750
$ws.Range("J46").Value = 4286.2856
$ws.Range("K46").Value = 2750
$ws.Range("L46").Value = 4286.2856
$ws.Range("M46").Value = -2562
$ws.Range("N46").Value = -4662.2856
$ws.Range("H100").Value = 6166.6665
$ws.Range("I100").Value = 4250
$ws.Range("K100").Value = 4250
$ws.Range("M100").Value = -3709
$ws.Range("H122").Value = 4999.25
$ws.Range("I122").Value = 4999.25
$ws.Range("K122").Value = 14997.75
$ws.Range("M122").Value = -12547.75
$ws.Range("H126").Value = 4336.875
$ws.Range("I126").Value = 3924.5
$ws.Range("K126").Value = 11773.5
$ws.Range("M126").Value = -9303.5
$ws.Range("H132").Value = 5500
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 602
$ws.Range("I14").Value = 602
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 602
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -434
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H18").Value = 24999
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H31").Value = 4900
$ws.Range("I31").Value = 4900
$ws.Range("K31").Value = 4900
$ws.Range("M31").Value = -4552
$ws.Range("H100").Value = 2596.8572
$ws.Range("I100").Value = 2449.1
$ws.Range("K100").Value = 4898.2
$ws.Range("M100").Value = -4357.2
$ws.Range("H122").Value = 5550.1
$ws.Range("I122").Value = 5187
$ws.Range("K122").Value = 15561
$ws.Range("M122").Value = -13111
$ws.Range("H126").Value = 2066.6
$ws.Range("I126").Value = 2066.6
$ws.Range("K126").Value = 6199.799999999999
$ws.Range("M126").Value = -3729.799999999999
